$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("NGIN")

# Replace the NGIN test data (row 2) - change the "1021" suffix/number to "1072"
$ws.Range("A2").Value  = "NGIN1072"
$ws.Range("C2").Value  = "NGIN1072"
$ws.Range("D2").Value  = "ngindomain1072.com"
$ws.Range("F2").Value  = "nginocn1072"
$ws.Range("G2").Value  = "testreference1072"
$ws.Range("H2").Value  = "ngincontact1072"
$ws.Range("J2").Value  = "ngin1072@test.com"
$ws.Range("Q2").Value  = "NGIN1072"
$ws.Range("T2").Value  = "NGINOrder_1072"
$ws.Range("U2").Value  = "NGINRFI_1072"
$ws.Range("V2").Value  = "NGINOrder_1072"
$ws.Range("W2").Value  = "NGINRFI_1072"
$ws.Range("Y2").Value  = "NGINService_1072"
$ws.Range("AJ2").Value = "NGINUser_1072"
$ws.Range("AK2").Value = "User_1072"
$ws.Range("AN2").Value = "NGINUser_1072@gmail.com"
$ws.Range("AP2").Value = "NGINOrder_1072"
$ws.Range("AQ2").Value = "NGINUseredit1072"
$ws.Range("AR2").Value = "Useredit1072"
$ws.Range("AU2").Value = "NGINUseredit_1072@gmail.com"
$ws.Range("AZ2").Value = "NGINOrderedit_1072"
$ws.Range("BA2").Value = "NGINRFIedit_1072"
$ws.Range("BB2").Value = "NGINOrder_1072"
$ws.Range("BC2").Value = "NGINRFI_1072"
$ws.Range("BI2").Value = "nginreseller1072@gmail.com"
$ws.Range("BQ2").Value = "nginreselleredit1072@gmail.com"
$ws.Range("BY2").Value = "Reseller1072"
$ws.Range("DP2").Value = "Reselleredit1072"
$ws.Range("ED2").Value = "AT-nginocn1072"
# Leading apostrophe preserves the existing "number stored as text" quote-prefix
# formatting on this cell (it was entered as text originally)
$ws.Range("EF2").Value = "'390201072891"

# Update the sheet view: scroll so column E is the leftmost visible column,
# and move the active selection to J5
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 5
$ws.Range("J5").Select()
